$d = $word.ActiveDocument

# 1. Fix the typo "Woopwoop!" -> "Woop!" in the first body paragraph.
$d.Content.Find.Execute("Woopwoop!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Woop!", 2)

# 2. Remove the now-obsolete "Ein hoffentlich finaler Git Action Test. Jetzt aber."
#    paragraph entirely (including its paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Ein hoffentlich finaler Git Action Test\. Jetzt aber\.") {
        $p.Range.Delete()
        break
    }
}
